# Update the "User stories (IO)" description in F15: it used to say the
# IO-version link also covers preceding versions; it now clarifies it
# covers the IO version ONLY (not preceding versions).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F15").Value = "Link naar open en gesloten user stories die op de betreffende IO versie maar niet de voorgaande versies betrekking hebben."

# Update the "Features (IO)" description in F19: it used to say the
# features.md contains all features for the version and preceding
# versions; it now says it contains only the features for that version.
$ws.Range("F19").Value = "Link naar een features.md indien aanwezig. Bevat alleen features voor de betrefende versie."

# Move the active selection to G25 (matches the saved workbook view state).
$ws.Range("G25").Select()
